$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force numeric-looking Price cells to remain stored as text (matches original text cell type)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply updated values
$ws.Range('D2').Value = '28.928.58'
$ws.Range('E2').Value = '  -2.09%  '
$ws.Range('D3').Value = '1.900.34'
$ws.Range('E3').Value = '  -3.94%  '
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '324.15'
$ws.Range('E5').Value = '  -1.23%  '
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').Value = '0.4594'
$ws.Range('E7').Value = '  -1.50%  '
$ws.Range('D8').Value = '0.3809'
$ws.Range('E8').Value = '  -2.72%  '
$ws.Range('D9').Value = '0.07710'
$ws.Range('E9').Value = '  -3.04%  '
$ws.Range('D10').Value = '0.9734'
$ws.Range('E10').Value = '  -2.09%  '
$ws.Range('D11').Value = '21.93'
$ws.Range('E11').Value = '  -4.10%  '
$ws.Range('D12').Value = '1.919.49'
$ws.Range('E12').Value = '  -3.83%  '
$ws.Range('D13').Value = '6.918'
$ws.Range('E13').Value = '  -3.81%  '
$ws.Range('E14').Value = '  -3.12%  '
$ws.Range('D15').Value = '0.07068'
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = '83.81'
$ws.Range('E17').Value = '  -4.47%  '
$ws.Range('D18').Value = '0.000009466'
$ws.Range('E18').Value = '  -5.07%  '
$ws.Range('D19').Value = '16.61'
$ws.Range('E19').Value = '  -4.02%  '
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('D21').Value = '28.905.87'
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('D22').Value = '5.270'
$ws.Range('E22').Value = '  -5.23%  '
$ws.Range('D23').Value = '10.84'
$ws.Range('E23').Value = '  -3.10%  '
$ws.Range('D24').Value = '2.099'
$ws.Range('E24').Value = '  -0.51%  '
$ws.Range('D25').Value = '158.15'
$ws.Range('E25').Value = '  -0.47%  '
$ws.Range('D26').Value = '19.02'
$ws.Range('E26').Value = '  -3.05%  '
$ws.Range('D27').Value = '5.610'
$ws.Range('E27').Value = '  -3.87%  '
$ws.Range('D28').Value = '117.48'
$ws.Range('E28').Value = '  -1.65%  '
$ws.Range('D29').Value = '1.841'
$ws.Range('E29').Value = '  -3.06%  '
$ws.Range('D30').Value = '0.09246'
$ws.Range('E30').Value = '  -1.90%  '
$ws.Range('D31').Value = '0.8550'
$ws.Range('E31').Value = '  -4.20%  '
$ws.Range('D32').Value = '5.073'
$ws.Range('E32').Value = '  -3.14%  '
$ws.Range('D33').Value = '1.235'
$ws.Range('E33').Value = '  -6.83%  '
$ws.Range('D34').Value = '2.951'
$ws.Range('E34').Value = '  -7.81%  '
$ws.Range('D35').Value = '0.05667'
$ws.Range('E35').Value = '  -2.51%  '
$ws.Range('E36').Value = '  -3.04%  '
$ws.Range('D37').Value = '1.005'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').Value = '0.02028'
$ws.Range('E38').Value = '  -3.38%  '
$ws.Range('D39').Value = '0.5469'
$ws.Range('E39').Value = '  -4.76%  '
$ws.Range('D40').Value = '7.376'
$ws.Range('E40').Value = '  -5.73%  '
$ws.Range('D41').Value = '0.1750'
$ws.Range('E41').Value = '  -3.01%  '
$ws.Range('D42').Value = '9.258'
$ws.Range('E42').Value = '  -4.37%  '
$ws.Range('D43').Value = '2.761'
$ws.Range('E43').Value = '  -0.78%  '
$ws.Range('D44').Value = '0.5146'
$ws.Range('E44').Value = '  -4.18%  '
$ws.Range('D45').Value = '11.12'
$ws.Range('E45').Value = '  -5.88%  '
$ws.Range('D46').Value = '0.06821'
$ws.Range('E46').Value = '  -1.77%  '
$ws.Range('D47').Value = '2.056'
$ws.Range('E47').Value = '  -5.03%  '
$ws.Range('D48').Value = '0.000002559'
$ws.Range('E48').Value = '  -17.66%  '
$ws.Range('D49').Value = '110.03'
$ws.Range('E49').Value = '  -3.60%  '
$ws.Range('E50').Value = '  -3.55%  '
$ws.Range('E51').Value = '  -0.15%  '
